$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Seller name"
$ws.Range("D3").Value = "Quantity sold"
$ws.Range("E3").Value = "Price (USD)"
$ws.Range("F3").Value = "Revenue (USD)"

$ws.Range("C4").Value = "BlueTech Goods"
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 250

$ws.Range("C5").Value = "EcoTrend"
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 600

$ws.Range("C6").Value = "GadgetPro"
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 22
$ws.Range("F6").Value = 330

$ws.Range("C7").Value = "HomeEssentials"
$ws.Range("D7").Value = 25
$ws.Range("E7").Value = 27
$ws.Range("F7").Value = 675

[void]$ws.Range("F10").Select()
